$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the number format of the date column (B2:B4) from a date format
# to plain text, then store the dates as literal text strings instead of
# numeric date serials.
$ws.Range("B2:B4").NumberFormat = "@"

$ws.Range("B2").Value = "10/1/2023"
$ws.Range("B3").Value = "10/2/2023"
$ws.Range("B4").Value = "10/3/2023"

# Update the active selection to match the authored workbook.
$ws.Range("D5").Select()
